# Adds new inventory/order query submissions to the bottom of Sheet1.
# Columns: B=Email, C=Name, D=Committee, E=Event, F=Item Code, G=Description, H=Quantity
# All source data in this sheet is free-text form input. Values that look
# numeric (item codes / quantities) need to be forced to text so they keep
# being stored as plain text, matching the rest of the sheet, instead of
# Excel auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Looks-Numeric($Text) {
    return $Text -match '^-?[0-9]+(\.[0-9]+)?$'
}

function Set-Cell($Row, $Col, $Text) {
    if ($Text -eq "") { return }
    $cell = $ws.Cells.Item($Row, $Col)
    if (Looks-Numeric $Text) {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $Text
}

function Set-Row($Row, $Email, $Name, $Committee, $Event, $ItemCode, $Description, $Quantity) {
    Set-Cell $Row 2 $Email
    Set-Cell $Row 3 $Name
    Set-Cell $Row 4 $Committee
    Set-Cell $Row 5 $Event
    Set-Cell $Row 6 $ItemCode
    Set-Cell $Row 7 $Description
    Set-Cell $Row 8 $Quantity
}

Set-Row 22 "shahzeb2000786@gmail.com" "Shahzeb Ahmed" "OPP" "coolio" "434" "meraj" "4"

Set-Row 23 "shahzeb2000786@gmail.com" "Shahzeb Ahmed" "" "" "" "" "55"

Set-Row 24 "" "fsfdsf" "" "" "" "" ""

Set-Row 25 "ahmeds85165@live.bucks.edu" "Shahzeb Ahmed" "technology committee" "Parade" "A19" "random" "444"

Set-Row 26 "shahzeb2000786@gmail.com" "Shahzeb Ahmed" "tech" "Parade" "A12" "Aluminum" "6"

Set-Row 27 "shahzeb2000786@gmail.com" "Shahzeb Ahmed" "OPP" "Game" "A15" "cups" "55"

Set-Row 28 "shahzeb2000786@gmail.com" "Shahzeb Ahmed" "Tech" "Parade" "A88" "Cups" "5"
